# Weekly update: a new Ajo (Chino / Primera) price record for
# "Vega Monumental Concepción" needs to be inserted as the new row 103,
# pushing all the subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 103 (shifts rows 103:224 down to 104:225,
# carrying formatting such as the date style on column D along with them).
$ws.Rows(103).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A103").Value = 11
$ws.Range("B103").Value = "Vega Monumental Concepción"
$ws.Range("C103").Value = "Bíobío"
$ws.Range("D103").Value = 44880
$ws.Range("E103").Value = 8
$ws.Range("F103").Value = 100112003
$ws.Range("G103").Value = "Ajo"
$ws.Range("H103").Value = "Chino"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 14000
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = 14500
$ws.Range("N103").Value = "$/caja 10 kilos"
$ws.Range("O103").Value = "China"
$ws.Range("P103").Value = 1450
$ws.Range("Q103").Value = 10
$ws.Range("R103").Value = "Hortaliza"
